$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up a handful of existing "0" text cells that become real
#        numeric zeros in this update (rows 73-75). ---
$ws.Cells.Item(73, 8).Value = 0   # H73
$ws.Cells.Item(74, 8).Value = 0   # H74
$ws.Cells.Item(75, 6).Value = 0   # F75
$ws.Cells.Item(75, 8).Value = 0   # H75

# --- 2. Append the new 2025-10-30 circular block (rows 87-103), scraped
#        from the updated PDF URL. This duplicates rows 70-86 verbatim
#        except for the B-column pdf_url, so build it from a small table. ---

# Force the new block to be plain text so Excel's COM layer doesn't
# "helpfully" reinterpret date-like / numeric-looking strings (dates in
# column A, thousands-separated numbers elsewhere) as real numbers.
$ws.Range("A87:I103").NumberFormat = "@"

$pdfUrl = "https://rashtriyametal.com/wp-content/uploads/2025/11/ListPrice30102025.pdf"
$circDate = "2025-10-30"

$rows = @(
    @{ C = "CHANDERIYA `nLEAD ZINC `nSMELTER";            D = "";  E = "330,100 331,600 330,600 329,600 328,100"; F = "";  G = "";             H = "";  I = "209,400" },
    @{ C = "HYDRO-1 UNIT";                                D = "";  E = "330,100 331,600 330,600 329,600 328,100"; F = "";  G = "";             H = "";  I = "209,400" },
    @{ C = "NEW HYDRO `nSMELTER `nCHANDERIYA";            D = "";  E = "330,100 331,600 330,600 329,600 328,100"; F = "";  G = "";             H = "";  I = "209,400" },
    @{ C = "ZINC SMELTER `nDEBRI";                        D = "0"; E = "0";                                       F = "";  G = "0  329,600";    H = "0"; I = "0" },
    @{ C = "Pantnagar `nMelting&Castin `ngPlant";         D = "";  E = "330,100 331,600 330,600 329,600";         F = "";  G = "";             H = "0"; I = "209,400" },
    @{ C = "RAJPURA DARIBA `nLEAD SMELTER";               D = "0"; E = "0";                                       F = "0"; G = "0";             H = "0"; I = "209,400" },
    @{ C = "Faridabad `nDepot";                           D = "";  E = "332,600 334,100 328,100 332,100 330,600"; F = "";  G = "";             H = "";  I = "211,900" },
    @{ C = "Panvel Depot";                                D = "";  E = "333,400 334,900 333,900 332,900 331,400"; F = "";  G = "";             H = "";  I = "212,300" },
    @{ C = "Pune Depot";                                  D = "";  E = "333,400 334,900 333,900 332,900 331,400"; F = "";  G = "";             H = "";  I = "212,700" },
    @{ C = "Baroda Depot";                                D = "";  E = "333,400 334,900 333,900 332,900 331,400"; F = "";  G = "";             H = "";  I = "212,700" },
    @{ C = "Raipur Depot";                                D = "";  E = "333,400 334,900 333,900 332,900 331,400"; F = "";  G = "";             H = "";  I = "212,700" },
    @{ C = "JAMSHEDPUR `nSTOCK POINT";                    D = "";  E = "331,100 332,600 331,600 330,600 329,100"; F = "";  G = "";             H = "";  I = "210,400" },
    @{ C = "";                                            D = "Kolkata Depot  331,100 332,600 331,600 330,600 329,100"; E = ""; F = ""; G = ""; H = "";  I = "210,400" },
    @{ C = "Bangalore `nDepot";                           D = "";  E = "331,100 332,600 331,600 330,600 329,100"; F = "";  G = "";             H = "";  I = "210,400" },
    @{ C = "Hyderabad `nDepot";                           D = "";  E = "331,100 332,600 331,600 330,600 329,100"; F = "";  G = "";             H = "";  I = "210,400" },
    @{ C = "";                                            D = "Chennai Depot  331,100 332,600 331,600 330,600 329,100"; E = ""; F = ""; G = ""; H = "";  I = "210,400" },
    @{ C = "Sindesar `nsmelter HZAPL";                    D = "";  E = "330,100 331,600";                         F = "";  G = "0  329,600 328,100"; H = ""; I = "209,400" }
)

$r = 87
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $circDate
    $ws.Cells.Item($r, 2).Value = $pdfUrl
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $r++
}
